$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

# Add the new location value to the next empty row in column A
$ws.Range("A10").Value = "Pita Kotte"

# Match the final selection recorded in the saved file
$ws.Range("F10").Select()
